$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for first new table (row 10)
$ws.Range("A10").Value = "flag[0]"
$ws.Range("B10").Value = "flag[1]"
$ws.Range("C10").Value = "lock[0]"
$ws.Range("D10").Value = "lock[1]"
$ws.Range("E10").Value = "thred 0"
$ws.Range("F10").Value = "thread 1"

# Row 11
$ws.Range("A11").Value = $false
$ws.Range("B11").Value = $false
$ws.Range("C11").Value = $false
$ws.Range("D11").Value = $false
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0

# Row 12
$ws.Range("A12").Value = $false
$ws.Range("B12").Value = $true
$ws.Range("C12").Value = $false
$ws.Range("D12").Value = $false
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 3

# Row 13
$ws.Range("A13").Value = $false
$ws.Range("B13").Value = $true
$ws.Range("C13").Value = $false
$ws.Range("D13").Value = $true
$ws.Range("E13").Value = 9
$ws.Range("F13").Value = 15

# Header row for second new table (row 16)
$ws.Range("A16").Value = "flag[0]"
$ws.Range("B16").Value = "flag[1]"
$ws.Range("C16").Value = "lock[0]"
$ws.Range("D16").Value = "lock[1]"
$ws.Range("E16").Value = "thred 0"
$ws.Range("F16").Value = "thread 1"

# Row 17
$ws.Range("A17").Value = $false
$ws.Range("B17").Value = $false
$ws.Range("C17").Value = $false
$ws.Range("D17").Value = $false
$ws.Range("E17").Value = 9
$ws.Range("F17").Value = 15

# Row 18
$ws.Range("A18").Value = $false
$ws.Range("B18").Value = $true
$ws.Range("C18").Value = $false
$ws.Range("D18").Value = $false
$ws.Range("E18").Value = 11
$ws.Range("F18").Value = 15

# Row 19
$ws.Range("A19").Value = $false
$ws.Range("B19").Value = $true
$ws.Range("C19").Value = $false
$ws.Range("D19").Value = $false
$ws.Range("E19").Value = 11
$ws.Range("F19").Value = 15

# Update sheet view: remove the topLeftCell freeze position and change selection
$ws.Range("A10:F13").Select() | Out-Null
